# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime" /
# "Latest HO Xliff Generate Date" timestamps for the ca6f5501-... file row
# (row 3) across the Overview, zh-cn and de-de sheets, reflecting a newly
# generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-30 04:48:18"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-30 04:48:13"
$wsZhCn.Range("K3").Value = "2016-08-30 04:48:28"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-30 04:48:18"
$wsDeDe.Range("K3").Value = "2016-08-30 04:48:35"
